$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data point was recorded for this market/product. It belongs
# right where row 115 currently sits (same category/variety as the row
# already there), so push the existing row 115 (and everything below it)
# down by one and populate the freed-up row with the new record.
$ws.Rows.Item(115).Insert()

$ws.Cells.Item(115, 1).Value = 5
$ws.Cells.Item(115, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(115, 3).Value = "Maule"
$ws.Cells.Item(115, 4).Value = 44601
$ws.Cells.Item(115, 5).Value = 7
$ws.Cells.Item(115, 6).Value = 100112003
$ws.Cells.Item(115, 7).Value = "Ajo"
$ws.Cells.Item(115, 8).Value = "Chino"
$ws.Cells.Item(115, 9).Value = "Primera"
$ws.Cells.Item(115, 10).Value = 200
$ws.Cells.Item(115, 11).Value = 20000
$ws.Cells.Item(115, 12).Value = 20000
$ws.Cells.Item(115, 13).Value = 20000
$ws.Cells.Item(115, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(115, 15).Value = "China"
$ws.Cells.Item(115, 16).Value = 2000
$ws.Cells.Item(115, 17).Value = 10
$ws.Cells.Item(115, 18).Value = "Hortaliza"
